# Loan RBI, Variable Instalments
# Insert a new (blank) column before column N on the "Repayment schedule"
# sheet - shifting the existing "Late" / "Original" / "Outstanding" columns
# one place to the right - and make "Repayment schedule" the active sheet
# (it was "Transactions" before).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a blank column at N; everything from N onward shifts right to O.
$ws.Columns("N").Insert()

# New column keeps the same on-screen width as the neighbouring columns
# (stored width 11, i.e. ColumnWidth 10.1666... in Excel's character units).
$ws.Columns("N").ColumnWidth = 10.166666666666666

# "Repayment schedule" becomes the selected/active sheet and tab.
$ws.Activate()
$ws.Range("S7").Select() | Out-Null
